$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text
$ws.Range("B1").Value = "Value (g)"

# Convert specific cells (numbers) to text values with "clean" decimal representation
$cells = @{
    "B7"  = "115.6"
    "B8"  = "74.8"
    "B9"  = "96.8"
    "B15" = "118.8"
    "B19" = "75.6"
    "B24" = "26.8"
    "B25" = "34.4"
    "B26" = "118.8"
    "B27" = "9.6"
    "B29" = "45.2"
    "B32" = "132.8"
}

foreach ($addr in $cells.Keys) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = [double]$cells[$addr]
    $c.ClearFormats()
}

# Update selection to E6
$ws.Range("E6").Select()
